# Bulk upload voters is now working!!!
# Rename the header row from "Full name"/"Email" to the lowercase,
# underscore-separated field names expected by the bulk-upload importer
# ("full_name"/"email"), leaving the voter rows untouched, and update the
# active cell selection to reflect where the user last clicked.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "full_name"
$ws.Range("B1").Value = "email"

$ws.Range("B12").Select()
